$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2453.2856
$ws.Range("I18").Value = 2304.7
$ws.Range("J18").Value = 2824.75
$ws.Range("K18").Value = 2304.7
$ws.Range("L18").Value = 2824.75
$ws.Range("M18").Value = -2020.7
$ws.Range("N18").Value = -3392.75
$ws.Range("H26").Value = 18000
$ws.Range("J26").Value = 18000
$ws.Range("L26").Value = 18000
$ws.Range("N26").Value = -18688
$ws.Range("H32").Value = 1168.0667
$ws.Range("J32").Value = 1092.4
$ws.Range("L32").Value = 1092.4
$ws.Range("N32").Value = -1744.4
$ws.Range("H40").Value = 1819
$ws.Range("I40").Value = 2301.5833
$ws.Range("J40").Value = 1405.3572
$ws.Range("K40").Value = 2301.5833
$ws.Range("L40").Value = 1405.3572
$ws.Range("M40").Value = -2126.5833
$ws.Range("N40").Value = -1755.3572
$ws.Range("H69").Value = 4309.375
$ws.Range("I69").Value = 5125
$ws.Range("K69").Value = 15375
$ws.Range("M69").Value = -14501
$ws.Range("H72").Value = 4309.375
$ws.Range("I72").Value = 5125
$ws.Range("K72").Value = 46125
$ws.Range("M72").Value = -41757
$ws.Range("H94").Value = 5111.1113
$ws.Range("H96").Value = 1651.0769
$ws.Range("I96").Value = 2545.2
$ws.Range("J96").Value = 1092.25
$ws.Range("K96").Value = 7635.599999999999
$ws.Range("L96").Value = 3276.75
$ws.Range("M96").Value = -6262.599999999999
$ws.Range("N96").Value = -6022.75
$ws.Range("H100").Value = 2235.4666
$ws.Range("I100").Value = 1893.2
$ws.Range("J100").Value = 2920
$ws.Range("K100").Value = 1893.2
$ws.Range("L100").Value = 2920
$ws.Range("M100").Value = -1352.2
$ws.Range("N100").Value = -4002
$ws.Range("H138").Value = 6580652
$ws.Range("I138").Value = 1285
$ws.Range("J138").Value = 16669014
$ws.Range("K138").Value = 3855
$ws.Range("L138").Value = 50007042
$ws.Range("M138").Value = 1285
$ws.Range("N138").Value = -50017322

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10297.804
$ws.Range("I32").Value = 11070.182
$ws.Range("J32").Value = 5442.857
$ws.Range("K32").Value = 11070.182
$ws.Range("L32").Value = 5442.857
$ws.Range("M32").Value = -10783.182
$ws.Range("N32").Value = -6016.857
$ws.Range("H37").Value = 12788
$ws.Range("I37").Value = 1538
$ws.Range("J37").Value = 17609.428
$ws.Range("K37").Value = 1538
$ws.Range("L37").Value = 17609.428
$ws.Range("M37").Value = -1265
$ws.Range("N37").Value = -18155.428

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12822087
$ws.Range("I86").Value = 1615.7084
$ws.Range("J86").Value = 33334840
$ws.Range("K86").Value = 1615.7084
$ws.Range("L86").Value = 33334840
$ws.Range("M86").Value = -492.7084
$ws.Range("N86").Value = -33337086
$ws.Range("H89").Value = 12822087
$ws.Range("I89").Value = 1615.7084
$ws.Range("J89").Value = 33334840
$ws.Range("K89").Value = 8078.541999999999
$ws.Range("L89").Value = 166674200
$ws.Range("M89").Value = -2462.541999999999
$ws.Range("N89").Value = -166685432
$ws.Range("H115").Value = 48684
$ws.Range("J115").Value = 48684
$ws.Range("L115").Value = 48684
$ws.Range("N115").Value = -51818

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 38992.6
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 38992.6
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 38992.6
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -39566.6
$ws.Range("H31").Value = 8777064
$ws.Range("I31").Value = 4773.486
$ws.Range("J31").Value = 111120450
$ws.Range("K31").Value = 4773.486
$ws.Range("L31").Value = 111120450
$ws.Range("M31").Value = -4478.486
$ws.Range("N31").Value = -111121040
$ws.Range("H34").Value = 8777064
$ws.Range("I34").Value = 4773.486
$ws.Range("J34").Value = 111120450
$ws.Range("K34").Value = 4773.486
$ws.Range("L34").Value = 111120450
$ws.Range("M34").Value = -4571.486
$ws.Range("N34").Value = -111120854
$ws.Range("H115").Value = 40290
$ws.Range("J115").Value = 40290
$ws.Range("L115").Value = 40290
$ws.Range("N115").Value = -42640

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2727.8
$ws.Range("I5").Value = 655.6
$ws.Range("K5").Value = 1966.8
$ws.Range("M5").Value = -1854.8
$ws.Range("H115").Value = 3065.4443
$ws.Range("I115").Value = 883
$ws.Range("J115").Value = 4156.6665
$ws.Range("K115").Value = 2649
$ws.Range("L115").Value = 12469.9995
$ws.Range("M115").Value = -1474
$ws.Range("N115").Value = -14819.9995
$ws.Range("H134").Value = 4522.381
$ws.Range("I134").Value = 2906.3635
$ws.Range("J134").Value = 6300
$ws.Range("K134").Value = 8719.0905
$ws.Range("L134").Value = 18900
$ws.Range("M134").Value = -3649.0905
$ws.Range("N134").Value = -29040
$ws.Range("H135").Value = 2727.8
$ws.Range("I135").Value = 655.6
$ws.Range("K135").Value = 5900.400000000001
$ws.Range("M135").Value = -3365.400000000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19059
$ws.Range("J15").Value = 19059
$ws.Range("L15").Value = 19059
$ws.Range("N15").Value = -19635
$ws.Range("H32").Value = 20000
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20592
$ws.Range("H81").Value = 19059
$ws.Range("J81").Value = 19059
$ws.Range("L81").Value = 19059
$ws.Range("N81").Value = -21055
$ws.Range("H84").Value = 19059
$ws.Range("J84").Value = 19059
$ws.Range("L84").Value = 57177
$ws.Range("N84").Value = -67161
$ws.Range("H102").Value = 1895.2444
$ws.Range("I102").Value = 2389.0356
$ws.Range("J102").Value = 1081.9412
$ws.Range("K102").Value = 2389.0356
$ws.Range("L102").Value = 1081.9412
$ws.Range("M102").Value = -767.0356000000002
$ws.Range("N102").Value = -4325.9412
$ws.Range("H103").Value = 7000
$ws.Range("J103").Value = 7000
$ws.Range("L103").Value = 7000
$ws.Range("N103").Value = -9344
$ws.Range("H111").Value = 23764
$ws.Range("J111").Value = 23764
$ws.Range("L111").Value = 23764
$ws.Range("N111").Value = -29898
$ws.Range("H113").Value = 1931
$ws.Range("I113").Value = 1805.5
$ws.Range("J113").Value = 2056.5
$ws.Range("K113").Value = 1805.5
$ws.Range("L113").Value = 2056.5
$ws.Range("M113").Value = 364.5
$ws.Range("N113").Value = -6396.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 426.55554
$ws.Range("I55").Value = 272.83334
$ws.Range("J55").Value = 503.41666
$ws.Range("K55").Value = 272.83334
$ws.Range("L55").Value = 503.41666
$ws.Range("M55").Value = -99.83334000000002
$ws.Range("N55").Value = -849.41666
$ws.Range("H61").Value = 1699.2307
$ws.Range("I61").Value = 1659
$ws.Range("J61").Value = 1833.3334
$ws.Range("K61").Value = 1659
$ws.Range("L61").Value = 1833.3334
$ws.Range("M61").Value = -1457
$ws.Range("N61").Value = -2237.3334
$ws.Range("H113").Value = 1699.2307
$ws.Range("I113").Value = 1659
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 1659
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = 511
$ws.Range("N113").Value = -6173.3334

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2753.3333
$ws.Range("I113").Value = 633.3333
$ws.Range("J113").Value = 3813.3333
$ws.Range("K113").Value = 1899.9999
$ws.Range("L113").Value = 11439.9999
$ws.Range("M113").Value = 270.0001
$ws.Range("N113").Value = -15779.9999
